$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.32821866666667
$ws.Range("H2").Value = 60.984656
$ws.Range("I2").Value = 0.004181898474048532
$ws.Range("J2").Value = 0.004181898474048532
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 630.0042926727823
$ws.Range("R2").Value = 5670.03863405504
$ws.Range("S2").Value = 0.001505116270618137
$ws.Range("T2").Value = 0.001505116270618137

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.32821866666667
$ws.Range("H3").Value = 60.984656
$ws.Range("I3").Value = 0.004181898474048532
$ws.Range("J3").Value = 0.004181898474048532
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 608.0834732668213
$ws.Range("R3").Value = 5472.751259401392
$ws.Range("S3").Value = 0.001452746179910946
$ws.Range("T3").Value = 0.001452746179910946

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.32821866666667
$ws.Range("H4").Value = 60.984656
$ws.Range("I4").Value = 0.004181898474048532
$ws.Range("J4").Value = 0.004181898474048532
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 389.289819018336
$ws.Range("R4").Value = 3503.608371165024
$ws.Range("S4").Value = 0.0009300356321457824
$ws.Range("T4").Value = 0.0009300356321457822

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.32821866666667
$ws.Range("H5").Value = 60.984656
$ws.Range("I5").Value = 0.004181898474048532
$ws.Range("J5").Value = 0.004181898474048532
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 123.0612625938987
$ws.Range("R5").Value = 1107.551363345088
$ws.Range("S5").Value = 0.0002940003913736668
$ws.Range("T5").Value = 0.0002940003913736667

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4809.896321333334
$ws.Range("H6").Value = 14429.688964
$ws.Range("I6").Value = 0.9894865072215304
$ws.Range("J6").Value = 0.9894865072215304
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 149066.4469641851
$ws.Range("R6").Value = 1341598.022677666
$ws.Range("S6").Value = 0.3561282634713128
$ws.Range("T6").Value = 0.3561282634713127

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4809.896321333334
$ws.Range("H7").Value = 14429.688964
$ws.Range("I7").Value = 0.9894865072215304
$ws.Range("J7").Value = 0.9894865072215304
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 143879.7225221545
$ws.Range("R7").Value = 1294917.50269939
$ws.Range("S7").Value = 0.3437368822700932
$ws.Range("T7").Value = 0.3437368822700932

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4809.896321333334
$ws.Range("H8").Value = 14429.688964
$ws.Range("I8").Value = 0.9894865072215304
$ws.Range("J8").Value = 0.9894865072215304
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 92110.56311093138
$ws.Range("R8").Value = 828995.0679983825
$ws.Range("S8").Value = 0.2200574009518191
$ws.Range("T8").Value = 0.2200574009518191

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4809.896321333334
$ws.Range("H9").Value = 14429.688964
$ws.Range("I9").Value = 0.9894865072215304
$ws.Range("J9").Value = 0.9894865072215304
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 29117.74631879674
$ws.Range("R9").Value = 262059.7168691707
$ws.Range("S9").Value = 0.06956396052830535
$ws.Range("T9").Value = 0.06956396052830534

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.69506
$ws.Range("H10").Value = 8.085180000000001
$ws.Range("I10").Value = 0.000554424737665286
$ws.Range("J10").Value = 0.000554424737665286
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 83.52425743013335
$ws.Range("R10").Value = 751.7183168712002
$ws.Range("S10").Value = 0.0001995442258274992
$ws.Range("T10").Value = 0.0001995442258274991

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.69506
$ws.Range("H11").Value = 8.085180000000001
$ws.Range("I11").Value = 0.000554424737665286
$ws.Range("J11").Value = 0.000554424737665286
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 80.61805475114001
$ws.Range("R11").Value = 725.5624927602601
$ws.Range("S11").Value = 0.0001926011414886457
$ws.Range("T11").Value = 0.0001926011414886456

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.69506
$ws.Range("H12").Value = 8.085180000000001
$ws.Range("I12").Value = 0.000554424737665286
$ws.Range("J12").Value = 0.000554424737665286
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 51.61098652308
$ws.Range("R12").Value = 464.4988787077201
$ws.Range("S12").Value = 0.0001233015972462391
$ws.Range("T12").Value = 0.000123301597246239

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.69506
$ws.Range("H13").Value = 8.085180000000001
$ws.Range("I13").Value = 0.000554424737665286
$ws.Range("J13").Value = 0.000554424737665286
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 16.31512784296
$ws.Range("R13").Value = 146.83615058664
$ws.Range("S13").Value = 0.00003897777310290219
$ws.Range("T13").Value = 0.00003897777310290219

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 28.08283533333333
$ws.Range("H14").Value = 84.24850599999999
$ws.Range("I14").Value = 0.005777169566755752
$ws.Range("J14").Value = 0.005777169566755752
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 870.3323739543378
$ws.Range("R14").Value = 7832.99136558904
$ws.Range("S14").Value = 0.002079273795622784
$ws.Range("T14").Value = 0.002079273795622783

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 28.08283533333333
$ws.Range("H15").Value = 84.24850599999999
$ws.Range("I15").Value = 0.005777169566755752
$ws.Range("J15").Value = 0.005777169566755752
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 840.0494076087047
$ws.Range("R15").Value = 7560.444668478342
$ws.Range("S15").Value = 0.002006926057838293
$ws.Range("T15").Value = 0.002006926057838293

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 28.08283533333333
$ws.Range("H16").Value = 84.24850599999999
$ws.Range("I16").Value = 0.005777169566755752
$ws.Range("J16").Value = 0.005777169566755752
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 537.792418691436
$ws.Range("R16").Value = 4840.131768222924
$ws.Range("S16").Value = 0.001284816832205264
$ws.Range("T16").Value = 0.001284816832205264

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 28.08283533333333
$ws.Range("H17").Value = 84.24850599999999
$ws.Range("I17").Value = 0.005777169566755752
$ws.Range("J17").Value = 0.005777169566755752
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 170.0055095827653
$ws.Range("R17").Value = 1530.049586244888
$ws.Range("S17").Value = 0.0004061528810894122
$ws.Range("T17").Value = 0.0004061528810894122
